$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup")

# Clear the shared-string values from the MQ_Experiment (column G) data
# cells that should be blanked out so autofill can repopulate them while
# keeping their existing cell style/formatting.
$blankRanges = @(
    "G3:G11",
    "G13:G21",
    "G23:G31",
    "G33:G41",
    "G43:G51",
    "G53:G61"
)

foreach ($rangeAddr in $blankRanges) {
    $ws.Range($rangeAddr).ClearContents() | Out-Null
}

# Move the active selection on the frozen (bottomLeft) pane from G50 to G51.
$ws.Activate() | Out-Null
$ws.Range("G51").Select() | Out-Null
